$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data was keyed in column-by-column (matches the order new strings
# land in the shared-string table: FMCG, ITC, Varun Beverages, ITC.NS, VBL.NS)

# Column A - S.No.
$ws.Cells.Item(15, 1).Value = 14
$ws.Cells.Item(16, 1).Value = 15
$ws.Cells.Item(17, 1).Value = 16

# Column B - Sector
$ws.Cells.Item(15, 2).Value = "FMCG"
$ws.Cells.Item(16, 2).Value = "FMCG"
$ws.Cells.Item(17, 2).Value = "Gold-ETF"

# Column C - Stock Name
$ws.Cells.Item(15, 3).Value = "ITC"
$ws.Cells.Item(16, 3).Value = "Varun Beverages"
$ws.Cells.Item(17, 3).Value = "ICICI Prudential Gold ETF"

# Column D - Stock Symbol
$ws.Cells.Item(15, 4).Value = "ITC.NS"
$ws.Cells.Item(16, 4).Value = "VBL.NS"
$ws.Cells.Item(17, 4).Value = "GOLDIETF.NS"

# Column E - Date
$ws.Cells.Item(15, 5).Value = "2025-03-17"
$ws.Cells.Item(16, 5).Value = "2025-03-17"
$ws.Cells.Item(17, 5).Value = "2025-03-17"

# Column F - Action
$ws.Cells.Item(15, 6).Value = "Buy"
$ws.Cells.Item(16, 6).Value = "Buy"
$ws.Cells.Item(17, 6).Value = "Buy"

# Column G - Units
$ws.Cells.Item(15, 7).Value = 1
$ws.Cells.Item(16, 7).Value = 1
$ws.Cells.Item(17, 7).Value = 1

# Column H - Price per Unit
$ws.Cells.Item(15, 8).Value = 408.9
$ws.Cells.Item(16, 8).Value = 504.7
$ws.Cells.Item(17, 8).Value = 75.74

# Copy style formatting from row 14 to the new rows 15-17 so that
# date / currency formats and alignment match the rest of the table.
[void]$ws.Range("A14:H14").Copy()
[void]$ws.Range("A15:H17").PasteSpecial(-4122) # xlPasteFormats
$excel.CutCopyMode = $false

# Update selection to match the final state (active cell H17)
[void]$ws.Range("H17").Select()
